$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'237.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.958"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.355"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05558"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.367"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'6.454"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "7KuCoinTokenKCS"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8007"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'1.036"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1399"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07313"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03287"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02866"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09231"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001666"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.254"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04775"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "'0.01166"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18OneONEBestin24h"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "'0.006258"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.005059"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.001053"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001499"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "UpBots"
$ws.Range("C24").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D24").Value = "'0.007485"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23UpBotsUBXT"
$ws.Range("D27").Value = "'0.1287"
$ws.Range("D27").Style = "Normal"
$ws.Range("D41").Value = "'0.006974"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003499"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1036"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008795"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005439"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6798"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.03170"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
